$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")
$ws.Activate()

$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 2
$ws.Range("J4").Value = 1
$ws.Range("J5").Value = 2
$ws.Range("J6").Value = 1
$ws.Range("J7").Value = 2

$ws.Range("N8").Select()
